# Week 16 logged + season sim performed from Week 17.
# Updates per-game log strings (YDS sheet play-by-play deltas, ST sheet
# return-yardage logs) and the season-total numeric tables (OFF, DEF, ST,
# TURNS, PEN) to reflect the newly logged week.

$wb = $excel.ActiveWorkbook

# ---- YDS: append Week 16 play-by-play deltas to the running logs ----
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + " 8 -2 19 3 1 1 9 1 1 1 3 12 2 2 2 8 2 2 2"
$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + " 10 18 20 3 11 16 5 17 3 68 20 17 9 11 15 3 52 1 19 3 18 9 3 2 5 14 13 0 12 4 22 3 10 20 4 52"
$ws.Range("C2").Value2 = $ws.Range("C2").Value2 + " 8 2 6 1 0 7 -2 3 -5 3 4 2 2 -7 1"
$ws.Range("C3").Value2 = $ws.Range("C3").Value2 + " 11 17 5 15 6 4 9 -1 16 12 10 -1 18 28 11 11 4 7 12 17 25 18 14 4 6 7 14 5"

# ---- OFF: season totals through Week 16/17 ----
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value2 = 203
$ws.Range("D2").Value2 = 8
$ws.Range("E2").Value2 = 10
$ws.Range("F2").Value2 = 51
$ws.Range("G2").Value2 = 52
$ws.Range("I2").Value2 = 5
$ws.Range("J2").Value2 = 24
$ws.Range("L2").Value2 = 283
$ws.Range("M2").Value2 = 207
$ws.Range("O2").Value2 = 18
$ws.Range("P2").Value2 = 11
$ws.Range("Q2").Value2 = 534
$ws.Range("B3").Value2 = 7
$ws.Range("C3").Value2 = 153
$ws.Range("D3").Value2 = 6
$ws.Range("E3").Value2 = 30
$ws.Range("F3").Value2 = 100
$ws.Range("G3").Value2 = 35
$ws.Range("H3").Value2 = 35
$ws.Range("I3").Value2 = 48
$ws.Range("J3").Value2 = 53
$ws.Range("N3").Value2 = 31

# ---- DEF: season totals through Week 16/17 ----
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value2 = 160
$ws.Range("D2").Value2 = 8
$ws.Range("F2").Value2 = 51
$ws.Range("G2").Value2 = 42
$ws.Range("H2").Value2 = 2
$ws.Range("J2").Value2 = 25
$ws.Range("L2").Value2 = 288
$ws.Range("M2").Value2 = 199
$ws.Range("O2").Value2 = 16
$ws.Range("Q2").Value2 = 492
$ws.Range("B3").Value2 = 11
$ws.Range("C3").Value2 = 200
$ws.Range("E3").Value2 = 38
$ws.Range("F3").Value2 = 106
$ws.Range("G3").Value2 = 38
$ws.Range("H3").Value2 = 41
$ws.Range("I3").Value2 = 57
$ws.Range("J3").Value2 = 50
$ws.Range("N3").Value2 = 21

# ---- ST: season totals + return-yardage logs ----
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value2 = 84
$ws.Range("F2").Value2 = 44
$ws.Range("G2").Value2 = 42
$ws.Range("J2").Value2 = 51
$ws.Range("K2").Value2 = 49
$ws.Range("N2").Value2 = 21
$ws.Range("B3").Value2 = 54
$ws.Range("B4").Value2 = $ws.Range("B4").Value2 + " 66 56 64 63"
$ws.Range("B5").Value2 = $ws.Range("B5").Value2 + " 22 15 24 0"
$ws.Range("B6").Value2 = $ws.Range("B6").Value2 + " 25 20 0"
$ws.Range("D5").Value2 = $ws.Range("D5").Value2 + " 6 10 0"

# ---- TURNS: season totals ----
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value2 = 5
$ws.Range("E2").Value2 = 10
$ws.Range("E3").Value2 = 8

# ---- PEN: season totals ----
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value2 = 11
$ws.Range("D2").Value2 = 6
$ws.Range("B3").Value2 = 16
